$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the previously-missing X10/Y10 cells on the existing row 10.
$ws.Range("X10").Value = -3.4100040000000149
$ws.Range("Y10").Value = "Down"

# Append a brand-new row 11 of trading data.
$ws.Range("A11").Value = 42654.886678240742
$ws.Range("B11").Value = 7
$ws.Range("C11").Value = "Buy"
$ws.Range("D11").Value = 24
$ws.Range("E11").Value = 34692
$ws.Range("F11").Value = 3963
$ws.Range("G11").Value = 60
$ws.Range("H11").Value = 37
$ws.Range("I11").Value = 89
$ws.Range("J11").Value = 10
$ws.Range("K11").Value = 27129
$ws.Range("L11").Value = 338
$ws.Range("M11").Value = 209
$ws.Range("N11").Value = 94
$ws.Range("O11").Value = 11
$ws.Range("P11").Value = "Bag"
$ws.Range("Q11").Value = 39.313912976930268
$ws.Range("R11").Value = 1.8
$ws.Range("S11").Value = 0.086400000000000005
$ws.Range("T11").Value = -0.0115
$ws.Range("U11").Value = 5.85
$ws.Range("V11").Value = "N/A"
$ws.Range("W11").Value = 0

# Copy the existing date/percentage number formats onto the new row so the
# style table stays in sync with the rest of the sheet (rather than minting
# brand-new numFmt entries).
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("S9:T9").Copy()
$ws.Range("S11:T11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
